$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update end dates for existing board members (serial date numbers)
$ws.Range("C11").Value = 42957   # Conlon: end date extended to 8/10/2017
$ws.Range("C26").Value = 43039   # McKeever: end date extended to 10/31/2017
$ws.Range("C29").Value = 42308   # Rodriguez: end date extended to 10/31/2015

# Add a new board member row (Lightfoot), matching formatting of the row above
$ws.Range("A33").Copy()
$ws.Paste($ws.Range("A34"))
$ws.Range("A34").Value = "Lightfoot"

$ws.Range("B33:C33").Copy()
$ws.Paste($ws.Range("B34:C34"))
$ws.Range("B34").Value = 42194   # start 7/9/2015
$ws.Range("C34").Value = 42957   # end 8/10/2017

$ws.Range("D29").Select()
